$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "opp_no_iter" sheet as the 2nd sheet (between
#    "Opposite Pins" and "Same Side Pins").
# ---------------------------------------------------------------
$wsOpposite = $wb.Worksheets.Item(1)
$wsSameSide = $wb.Worksheets.Item(2)
$wsLattice  = $wb.Worksheets.Item(3)

$ws = $wb.Worksheets.Add($wsSameSide)
$ws.Name = "opp_no_iter"


# ---------------------------------------------------------------
# 2. Populate "opp_no_iter" with the converged (non-iterative)
#    opposite-pins results table.
# ---------------------------------------------------------------
# --- Header row 1 ---
$ws.Range("A1").Value = "Area"
$ws.Range("B1").Value = "N = 4"
$ws.Range("C1").Value = "N = 9"
$ws.Range("D1").Value = "N = 16"
$ws.Range("E1").Value = "N = 25"
$ws.Range("F1").Value = "N = 36"
$ws.Range("G1").Value = "N = 49"
$ws.Range("H1").Value = "N = 64"
$ws.Range("I1").Value = "N = 81"
$ws.Range("J1").Value = "N = 100"

# --- Column A (rows 2-19) ---
$ws.Range("A2").Value = "Coarse"
$ws.Range("A3").Value = 1.8
$ws.Range("A4").Value = 1.6
$ws.Range("A5").Value = 1.4
$ws.Range("A6").Value = 1.2
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 0.8
$ws.Range("A9").Value = 0.6
$ws.Range("A10").Value = 0.4
$ws.Range("A11").Value = 0.2
$ws.Range("A12").Value = 0.1
$ws.Range("A13").Value = 0.08
$ws.Range("A14").Value = 0.06
$ws.Range("A15").Value = 0.05
$ws.Range("A16").Value = 0.04
$ws.Range("A17").Value = 0.03
$ws.Range("A18").Value = 0.02
$ws.Range("A19").Value = 0.01

# --- Data columns B..J rows 2-19 ---
$ws.Range("B2").Value = 1.9487
$ws.Range("C2").Value = 4.1204
$ws.Range("D2").Value = 6.7555
$ws.Range("E2").Value = 9.5959
$ws.Range("F2").Value = 12.436
$ws.Range("G2").Value = 14.213
$ws.Range("H2").Value = 16.444
$ws.Range("I2").Value = 8.5977
$ws.Range("J2").Value = 6.77083
$ws.Range("B3").Value = 1.4575
$ws.Range("C3").Value = 2.3248
$ws.Range("D3").Value = 4.1121
$ws.Range("E3").Value = 4.6428
$ws.Range("F3").Value = 7.8387
$ws.Range("G3").Value = 8.6069
$ws.Range("H3").Value = 24.765
$ws.Range("I3").Value = 6.1363
$ws.Range("J3").Value = 4.57746
$ws.Range("B4").Value = 1.4227
$ws.Range("C4").Value = 2.2133
$ws.Range("D4").Value = 4.2027
$ws.Range("E4").Value = 4.6428
$ws.Range("F4").Value = 6.8571
$ws.Range("G4").Value = 8.518
$ws.Range("H4").Value = 24.715
$ws.Range("I4").Value = 5.9384
$ws.Range("J4").Value = 4.57746
$ws.Range("B5").Value = 1.3162
$ws.Range("C5").Value = 2.0536
$ws.Range("D5").Value = 2.9779
$ws.Range("E5").Value = 4.6428
$ws.Range("F5").Value = 6.2307
$ws.Range("G5").Value = 8.5771
$ws.Range("H5").Value = 19.983
$ws.Range("I5").Value = 5.9037
$ws.Range("J5").Value = 4.51389
$ws.Range("B6").Value = 1.2976
$ws.Range("C6").Value = 1.9524
$ws.Range("D6").Value = 3.0216
$ws.Range("E6").Value = 4.9295
$ws.Range("F6").Value = 4.5051
$ws.Range("G6").Value = 7.2459
$ws.Range("H6").Value = 19.973
$ws.Range("I6").Value = 4.2993
$ws.Range("J6").Value = 4.51389
$ws.Range("B7").Value = 1.3477
$ws.Range("C7").Value = 1.7544
$ws.Range("D7").Value = 2.9029
$ws.Range("E7").Value = 4.9295
$ws.Range("F7").Value = 4.5205
$ws.Range("G7").Value = 6.0226
$ws.Range("H7").Value = 20.013
$ws.Range("I7").Value = 4.6188
$ws.Range("J7").Value = 4.51389
$ws.Range("B8").Value = 1.2567
$ws.Range("C8").Value = 1.6526
$ws.Range("D8").Value = 2.9497
$ws.Range("E8").Value = 3.3112
$ws.Range("F8").Value = 4.4511
$ws.Range("G8").Value = 4.4001
$ws.Range("H8").Value = 19.739
$ws.Range("I8").Value = 4.5849
$ws.Range("J8").Value = 2.91667
$ws.Range("B9").Value = 1.1423
$ws.Range("C9").Value = 1.4516
$ws.Range("D9").Value = 2.0545
$ws.Range("E9").Value = 3.006
$ws.Range("F9").Value = 3.5497
$ws.Range("G9").Value = 4.2241
$ws.Range("H9").Value = 14.28
$ws.Range("I9").Value = 2.8738
$ws.Range("J9").Value = 3.09917
$ws.Range("B10").Value = 1.0945
$ws.Range("C10").Value = 1.3529
$ws.Range("D10").Value = 1.7868
$ws.Range("E10").Value = 2.0214
$ws.Range("F10").Value = 2.7379
$ws.Range("G10").Value = 3.3303
$ws.Range("H10").Value = 14.093
$ws.Range("I10").Value = 2.8027
$ws.Range("J10").Value = 2.06311
$ws.Range("B11").Value = 1.0457
$ws.Range("C11").Value = 1.1361
$ws.Range("D11").Value = 1.3358
$ws.Range("E11").Value = 1.5452
$ws.Range("F11").Value = 1.6476
$ws.Range("G11").Value = 2.0494
$ws.Range("H11").Value = 8.7796
$ws.Range("I11").Value = 1.8174
$ws.Range("J11").Value = 1.44558
$ws.Range("B12").Value = 1.0197
$ws.Range("C12").Value = 1.0426
$ws.Range("D12").Value = 1.1093
$ws.Range("E12").Value = 1.1698
$ws.Range("F12").Value = 1.2869
$ws.Range("G12").Value = 1.3647
$ws.Range("H12").Value = 4.4273
$ws.Range("I12").Value = 1.4091
$ws.Range("J12").Value = 1.24494
$ws.Range("B13").Value = 1.0112
$ws.Range("C13").Value = 1.0291
$ws.Range("D13").Value = 1.0943
$ws.Range("E13").Value = 1.19
$ws.Range("F13").Value = 1.2086
$ws.Range("G13").Value = 1.2899
$ws.Range("H13").Value = 3.3851
$ws.Range("I13").Value = 1.3178
$ws.Range("J13").Value = 1.17602
$ws.Range("B14").Value = 1.0052
$ws.Range("C14").Value = 1.0305
$ws.Range("D14").Value = 1.0373
$ws.Range("E14").Value = 1.1048
$ws.Range("F14").Value = 1.0866
$ws.Range("G14").Value = 1.1981
$ws.Range("H14").Value = 2.9331
$ws.Range("I14").Value = 1.2829
$ws.Range("J14").Value = 1.06347
$ws.Range("B15").Value = 1.0206
$ws.Range("C15").Value = 1.0216
$ws.Range("D15").Value = 1.0575
$ws.Range("E15").Value = 1.0916
$ws.Range("F15").Value = 1.0792
$ws.Range("G15").Value = 1.1146
$ws.Range("H15").Value = 2.6072
$ws.Range("I15").Value = 1.2175
$ws.Range("J15").Value = 1.08662
$ws.Range("B16").Value = 1.004
$ws.Range("C16").Value = 1.0134
$ws.Range("D16").Value = 1.0016
$ws.Range("E16").Value = 1.0608
$ws.Range("F16").Value = 1.0735
$ws.Range("G16").Value = 1.0731
$ws.Range("H16").Value = 2.2022
$ws.Range("I16").Value = 1.1711
$ws.Range("J16").Value = 1.10996
$ws.Range("B17").Value = 1.0029
$ws.Range("C17").Value = 1.0164
$ws.Range("D17").Value = 1.0208
$ws.Range("E17").Value = 1.0499
$ws.Range("F17").Value = 1.0651
$ws.Range("G17").Value = 1.0479
$ws.Range("H17").Value = 1.9293
$ws.Range("I17").Value = 1.1315
$ws.Range("J17").Value = 1.03406
$ws.Range("B18").Value = 1.0041
$ws.Range("C18").Value = 1.0082
$ws.Range("D18").Value = 1.0099
$ws.Range("E18").Value = 1.0337
$ws.Range("F18").Value = 1.024
$ws.Range("G18").Value = 1.0275
$ws.Range("H18").Value = 1.5743
$ws.Range("I18").Value = 1.0798
$ws.Range("J18").Value = 1.04855
$ws.Range("B19").Value = 1.0026
$ws.Range("C19").Value = 1.0101
$ws.Range("D19").Value = 1.0082
$ws.Range("E19").Value = 1.0088
$ws.Range("F19").Value = 1.0391
$ws.Range("G19").Value = 1.018
$ws.Range("H19").Value = 1.2765
$ws.Range("I19").Value = 1.0429
$ws.Range("J19").Value = 1.00967
# ---------------------------------------------------------------
# 3. Styling to match "Opposite Pins" / "Same Side Pins" headers:
#    - A1 and column A (A2:A19): bold, centered
#    - B1:J1: bold, centered horizontally + vertically
#    - J2:J19: centered (not bold)
# ---------------------------------------------------------------
$ws.Range("A1:J1").Font.Bold = $true
$ws.Range("A1:J1").HorizontalAlignment = -4108
$ws.Range("A2:A19").Font.Bold = $true
$ws.Range("A2:A19").HorizontalAlignment = -4108
$ws.Range("B1:J1").VerticalAlignment = -4108
$ws.Range("J2:J19").HorizontalAlignment = -4108


# ---------------------------------------------------------------
# 4. Sheet view / selection tweaks.
# ---------------------------------------------------------------
# New sheet: select the whole table.
$ws.Range("A1:J19").Select()

# "Opposite Pins": selection moves to A2:J20, no frozen/scrolled topLeftCell.
$wsOpposite.Range("A2:J20").Select()

# "Lattice": selection becomes A1:J19 (no longer the tab-selected sheet).
$wsLattice.Range("A1:J19").Select()

# ---------------------------------------------------------------
# 5. Make "opp_no_iter" the active tab/sheet (matches activeTab=1).
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("A1:J19").Select()
